$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: text values imported from the properties file
$ws.Range("A7").Value = "fgfdg"
$ws.Range("B7").Value = "tertret"

# Row 8: numeric values imported from the properties file
$ws.Range("A8").Value = 5454545
$ws.Range("B8").Font.Size = 14
$ws.Range("B8").HorizontalAlignment = -4131
$ws.Range("B8").VerticalAlignment = -4108
$ws.Range("B8").Value = 321
$ws.Rows.Item(8).RowHeight = 18

# Selection moves to A4
$ws.Range("A4").Select() | Out-Null
